# Auto-generated Excel COM-interop script applying the Ramuh_Profits market-data refresh.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for specific leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, matching a scheduled price-data sync.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1383.4445
$ws.Range("I43").Value = 1487.25
$ws.Range("J43").Value = 1300.4
$ws.Range("K43").Value = 1487.25
$ws.Range("L43").Value = 1300.4
$ws.Range("M43").Value = -1418.25
$ws.Range("N43").Value = -1438.4

$ws.Range("H55").Value = 166
$ws.Range("J55").Value = 140.66667
$ws.Range("L55").Value = 140.66667
$ws.Range("N55").Value = -568.6666700000001

$ws.Range("H86").Value = 4633.8335
$ws.Range("I86").Value = 4760.6
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 4760.6
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -3637.6
$ws.Range("N86").Value = -6246

$ws.Range("H89").Value = 4633.8335
$ws.Range("I89").Value = 4760.6
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 23803
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -18187
$ws.Range("N89").Value = -31232

$ws.Range("H100").Value = 85992.5
$ws.Range("J100").Value = 3190
$ws.Range("L100").Value = 3190
$ws.Range("N100").Value = -4272

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H138").Value = 4334.3184
$ws.Range("I138").Value = 1424.7407
$ws.Range("J138").Value = 8955.412
$ws.Range("K138").Value = 4274.2221
$ws.Range("L138").Value = 26866.236
$ws.Range("M138").Value = 865.7779
$ws.Range("N138").Value = -37146.236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 27544.75
$ws.Range("I6").Value = 27544.75
$ws.Range("K6").Value = 27544.75
$ws.Range("M6").Value = -27371.75

$ws.Range("H97").Value = 1034.6538
$ws.Range("I97").Value = 846.8946999999999
$ws.Range("J97").Value = 1544.2858
$ws.Range("K97").Value = 846.8946999999999
$ws.Range("L97").Value = 1544.2858
$ws.Range("M97").Value = -350.8946999999999
$ws.Range("N97").Value = -2536.2858

$ws.Range("H102").Value = 1489.7142
$ws.Range("I102").Value = 1107.5
$ws.Range("K102").Value = 1107.5
$ws.Range("M102").Value = 514.5

$ws.Range("H106").Value = 45120
$ws.Range("J106").Value = 45120
$ws.Range("L106").Value = 45120
$ws.Range("N106").Value = -47644

$ws.Range("H132").Value = 4481.3955
$ws.Range("I132").Value = 5611.923
$ws.Range("J132").Value = 2752.353
$ws.Range("K132").Value = 16835.769
$ws.Range("L132").Value = 8257.059000000001
$ws.Range("M132").Value = -14305.769
$ws.Range("N132").Value = -13317.059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4093.1304
$ws.Range("I86").Value = 3793.889
$ws.Range("J86").Value = 4285.5
$ws.Range("K86").Value = 3793.889
$ws.Range("L86").Value = 4285.5
$ws.Range("M86").Value = -2670.889
$ws.Range("N86").Value = -6531.5

$ws.Range("H89").Value = 4093.1304
$ws.Range("I89").Value = 3793.889
$ws.Range("J89").Value = 4285.5
$ws.Range("K89").Value = 18969.445
$ws.Range("L89").Value = 21427.5
$ws.Range("M89").Value = -13353.445
$ws.Range("N89").Value = -32659.5

$ws.Range("H94").Value = 1575.7142
$ws.Range("I94").Value = 1210
$ws.Range("J94").Value = 2490
$ws.Range("K94").Value = 1210
$ws.Range("L94").Value = 2490
$ws.Range("M94").Value = -759
$ws.Range("N94").Value = -3392

$ws.Range("H99").Value = 2550
$ws.Range("I99").Value = 2075
$ws.Range("K99").Value = 2075
$ws.Range("M99").Value = -577

$ws.Range("H106").Value = 19394.334
$ws.Range("I106").Value = 10000
$ws.Range("J106").Value = 20065.357
$ws.Range("K106").Value = 10000
$ws.Range("L106").Value = 20065.357
$ws.Range("M106").Value = -8738
$ws.Range("N106").Value = -22589.357

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 47574.332
$ws.Range("I135").Value = 30709
$ws.Range("J135").Value = 49682.5
$ws.Range("K135").Value = 30709
$ws.Range("L135").Value = 49682.5
$ws.Range("M135").Value = -25639
$ws.Range("N135").Value = -59822.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2108

$ws.Range("H64").Value = 2000

$ws.Range("H66").Value = 2108

$ws.Range("H67").Value = 2000

$ws.Range("H131").Value = 566.0909
$ws.Range("I131").Value = 264.75
$ws.Range("J131").Value = 899.4894
$ws.Range("K131").Value = 794.25
$ws.Range("L131").Value = 2698.4682
$ws.Range("M131").Value = 4245.75
$ws.Range("N131").Value = -12778.4682

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 23.52381
$ws.Range("I2").Value = 18.615385
$ws.Range("K2").Value = 18.615385
$ws.Range("M2").Value = 94.384615

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H80").Value = 4801.6665
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 5202.5
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 5202.5
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -7198.5

$ws.Range("H83").Value = 4801.6665
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 5202.5
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 26012.5
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -35996.5

$ws.Range("H97").Value = 607.4
$ws.Range("I97").Value = 544
$ws.Range("J97").Value = 966.6667
$ws.Range("K97").Value = 544
$ws.Range("L97").Value = 966.6667
$ws.Range("M97").Value = -48
$ws.Range("N97").Value = -1958.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 48000
$ws.Range("J75").Value = 48000
$ws.Range("L75").Value = 48000
$ws.Range("N75").Value = -49872

$ws.Range("H78").Value = 48000
$ws.Range("J78").Value = 48000
$ws.Range("L78").Value = 144000
$ws.Range("N78").Value = -153360

$ws.Range("H82").Value = 1334
$ws.Range("I82").Value = 1334
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1334
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -973
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 1334
$ws.Range("I85").Value = 1334
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1334
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -86
$ws.Range("N85").ClearContents()

$ws.Range("H100").Value = 1176.5
$ws.Range("I100").Value = 1176.5
$ws.Range("K100").Value = 1176.5
$ws.Range("M100").Value = -635.5

$ws.Range("H132").Value = 6210.9707
$ws.Range("I132").Value = 7605.55
$ws.Range("K132").Value = 22816.65
$ws.Range("M132").Value = -20286.65

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2360.8667
$ws.Range("I96").Value = 1628.1666
$ws.Range("K96").Value = 1628.1666
$ws.Range("M96").Value = -255.1666

Write-Host "Applied Ramuh_Profits market-data sync."